# Tesla + Performance Updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New portfolio entry: Tesla (TSLA), bought 10/23/24
$ws.Range("A7").Value = "Tesla"
$ws.Range("B7").Value = "TSLA"
$ws.Range("C7").Value = 216.58
$ws.Range("E7").Value = 205.22
$ws.Range("F7").Value = 239.1
$ws.Range("G7").Value = 280.12

# L7 already carries style s="9" (numFmtId "#,##0"). Typing a slash-date
# literal straight into Value auto-detects as a real date and rewrites the
# cell's style (new numFmt + quotePrefix + lost border). Temporarily switch
# to a text format, assign the literal, then restore the original number
# format so the cell keeps its original style index instead of minting a
# new one.
$L7 = $ws.Range("L7")
$L7.NumberFormat = "@"
$L7.Value = "10/23/24"
$L7.NumberFormat = "#,##0"

$ws.Range("M7").Value = 216.29
$ws.Range("N7").Value = "#"
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = 0

# Leave the active selection on P13, matching the author's last interaction
$ws.Range("P13").Select()
